$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Octubre de 2020 a las 11:36"

# Row 24 (Pakistan)
$ws.Range("B24").Value = 315260
$ws.Range("C24").Value = 644
$ws.Range("D24").Value = 299836
$ws.Range("E24").Value = 8907
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 6517

# Row 25 (Indonesia)
$ws.Range("B25").Value = 307120
$ws.Range("C25").Value = 3622
$ws.Range("D25").Value = 232593
$ws.Range("E25").Value = 63274
$ws.Range("G25").Value = 102
$ws.Range("H25").Value = 11253

# Rows 42/43: Oman and Polonia swap positions (Polonia now ranks above Oman)
# Row 42 becomes Polonia with fresh data; Row 43 becomes Oman with Polonia's old row-42 data
$ws.Range("A42").Value = "Polonia"
$ws.Range("B42").Value = 102080
$ws.Range("C42").Value = 2006
$ws.Range("D42").Value = 73552
$ws.Range("E42").Value = 25869
$ws.Range("G42").Value = 29
$ws.Range("H42").Value = 2659

$ws.Range("A43").Value = "Oman"
$ws.Range("B43").Value = 101270
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 90296
$ws.Range("E43").Value = 9997
$ws.Range("H43").Value = 977

# Row 60 (Singapur)
$ws.Range("B60").Value = 57819
$ws.Range("C60").Value = 7
$ws.Range("E60").Value = 217

# Row 65 (Austria)
$ws.Range("B65").Value = 48896
$ws.Range("C65").Value = 750
$ws.Range("D65").Value = 39058
$ws.Range("E65").Value = 9020
$ws.Range("G65").Value = 5
$ws.Range("H65").Value = 818

# Row 113
$ws.Range("E113").Value = 4024
$ws.Range("G113").Value = 3
$ws.Range("H113").Value = 53

# Row 117 (Georgia)
$ws.Range("B117").Value = 6673
$ws.Range("C117").Value = 175
$ws.Range("E117").Value = 2252
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 156

# Row 125 (Lituania)
$ws.Range("B125").Value = 5285
$ws.Range("C125").Value = 100
$ws.Range("D125").Value = 2497
$ws.Range("E125").Value = 2694

# Rows 215/216: Islas Malvinas and Montserrat swap positions
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
